# Add a "formula" column (E) that exercises error handling in formula cells.
# Mirrors columns A-D (double/boolean/timestamp/string) with a new
# CONCAT(...)-based formula column, including rows that raise #DIV/0! and
# #NAME? errors - matching the commit "handle errors in formula cells".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: E1 = "formula" (reuse D1's header style) ---
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "formula"

# --- Row 2: simple success case -> "A1" ---
$ws.Range("D2").Copy($ws.Range("E2"))
$ws.Range("E2").Formula = '=_xlfn.CONCAT("A", 3/3)'

# --- Row 3: simple success case -> "A3" ---
$ws.Range("D3").Copy($ws.Range("E3"))
$ws.Range("E3").Formula = '=_xlfn.CONCAT("A", 3/1)'

# --- Row 4: division by zero -> #DIV/0! ---
$ws.Range("D4").Copy($ws.Range("E4"))
$ws.Range("E4").Formula = '=_xlfn.CONCAT("A", 3/0)'

# --- Row 5: legacy array-entered NA (no parens) -> #NAME? ---
$ws.Range("D5").Copy($ws.Range("E5"))
$ws.Range("E5").FormulaArray = '=NA'
